$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3361.918
$ws.Cells.Item(76, 9).Value = 3126.8545
$ws.Cells.Item(76, 10).Value = 5516.6665
$ws.Cells.Item(76, 11).Value = 3126.8545
$ws.Cells.Item(76, 12).Value = 5516.6665
$ws.Cells.Item(76, 13).Value = -2811.8545
$ws.Cells.Item(76, 14).Value = -6146.6665
$ws.Cells.Item(79, 8).Value = 3361.918
$ws.Cells.Item(79, 9).Value = 3126.8545
$ws.Cells.Item(79, 10).Value = 5516.6665
$ws.Cells.Item(79, 11).Value = 3126.8545
$ws.Cells.Item(79, 12).Value = 5516.6665
$ws.Cells.Item(79, 13).Value = -2034.8545
$ws.Cells.Item(79, 14).Value = -7700.6665
$ws.Cells.Item(87, 8).Value = 18112.846
$ws.Cells.Item(87, 10).Value = 18112.846
$ws.Cells.Item(87, 12).Value = 18112.846
$ws.Cells.Item(87, 14).Value = -20608.846
$ws.Cells.Item(90, 8).Value = 18112.846
$ws.Cells.Item(90, 10).Value = 18112.846
$ws.Cells.Item(90, 12).Value = 54338.538
$ws.Cells.Item(90, 14).Value = -66818.538
$ws.Cells.Item(137, 8).Value = 1774.4482
$ws.Cells.Item(137, 9).Value = 1794.25
$ws.Cells.Item(137, 10).Value = 1679.4
$ws.Cells.Item(137, 11).Value = 5382.75
$ws.Cells.Item(137, 12).Value = 5038.200000000001
$ws.Cells.Item(137, 13).Value = -2832.75
$ws.Cells.Item(137, 14).Value = -10138.2
$ws.Cells.Item(141, 8).Value = 1817
$ws.Cells.Item(141, 9).Value = 1817
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 5451
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = -271
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4150.604
$ws.Cells.Item(61, 9).Value = 5560.8213
$ws.Cells.Item(61, 10).Value = 2571.16
$ws.Cells.Item(61, 11).Value = 5560.8213
$ws.Cells.Item(61, 12).Value = 2571.16
$ws.Cells.Item(61, 13).Value = -5348.8213
$ws.Cells.Item(61, 14).Value = -2995.16
$ws.Cells.Item(88, 8).Value = 7080.1816
$ws.Cells.Item(88, 9).Value = 11215.272
$ws.Cells.Item(88, 10).Value = 2945.0908
$ws.Cells.Item(88, 11).Value = 11215.272
$ws.Cells.Item(88, 12).Value = 2945.0908
$ws.Cells.Item(88, 13).Value = -10809.272
$ws.Cells.Item(88, 14).Value = -3757.0908
$ws.Cells.Item(91, 8).Value = 7080.1816
$ws.Cells.Item(91, 9).Value = 11215.272
$ws.Cells.Item(91, 10).Value = 2945.0908
$ws.Cells.Item(91, 11).Value = 11215.272
$ws.Cells.Item(91, 12).Value = 2945.0908
$ws.Cells.Item(91, 13).Value = -9811.272000000001
$ws.Cells.Item(91, 14).Value = -5753.0908
$ws.Cells.Item(110, 8).Value = 86835
$ws.Cells.Item(110, 10).Value = 2004
$ws.Cells.Item(110, 12).Value = 2004
$ws.Cells.Item(110, 14).Value = -6094
$ws.Cells.Item(132, 8).Value = 3012.6177
$ws.Cells.Item(132, 9).Value = 1458.6666
$ws.Cells.Item(132, 10).Value = 3860.2273
$ws.Cells.Item(132, 11).Value = 4375.9998
$ws.Cells.Item(132, 12).Value = 11580.6819
$ws.Cells.Item(132, 13).Value = -1845.9998
$ws.Cells.Item(132, 14).Value = -16640.6819
$ws.Cells.Item(136, 8).Value = 4150.604
$ws.Cells.Item(136, 9).Value = 5560.8213
$ws.Cells.Item(136, 10).Value = 2571.16
$ws.Cells.Item(136, 11).Value = 16682.4639
$ws.Cells.Item(136, 12).Value = 7713.48
$ws.Cells.Item(136, 13).Value = -14132.4639
$ws.Cells.Item(136, 14).Value = -12813.48

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(109, 8).Value = 63155
$ws.Cells.Item(109, 10).Value = 63155
$ws.Cells.Item(109, 12).Value = 63155
$ws.Cells.Item(109, 14).Value = -65929
$ws.Cells.Item(134, 8).Value = 2091.1943
$ws.Cells.Item(134, 9).Value = 1371.8846
$ws.Cells.Item(134, 10).Value = 3961.4
$ws.Cells.Item(134, 11).Value = 4115.6538
$ws.Cells.Item(134, 12).Value = 11884.2
$ws.Cells.Item(134, 13).Value = -1580.6538
$ws.Cells.Item(134, 14).Value = -16954.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 3546.818
$ws.Cells.Item(86, 9).Value = 4484.5
$ws.Cells.Item(86, 10).Value = 2421.6
$ws.Cells.Item(86, 11).Value = 4484.5
$ws.Cells.Item(86, 12).Value = 2421.6
$ws.Cells.Item(86, 13).Value = -3361.5
$ws.Cells.Item(86, 14).Value = -4667.6
$ws.Cells.Item(89, 8).Value = 3546.818
$ws.Cells.Item(89, 9).Value = 4484.5
$ws.Cells.Item(89, 10).Value = 2421.6
$ws.Cells.Item(89, 11).Value = 22422.5
$ws.Cells.Item(89, 12).Value = 12108
$ws.Cells.Item(89, 13).Value = -16806.5
$ws.Cells.Item(89, 14).Value = -23340
$ws.Cells.Item(132, 8).Value = 3720.125
$ws.Cells.Item(132, 9).Value = 2416.3333
$ws.Cells.Item(132, 10).Value = 4502.4
$ws.Cells.Item(132, 11).Value = 7248.999899999999
$ws.Cells.Item(132, 12).Value = 13507.2
$ws.Cells.Item(132, 13).Value = -4718.999899999999
$ws.Cells.Item(132, 14).Value = -18567.2
$ws.Cells.Item(134, 8).Value = 6633.7896
$ws.Cells.Item(134, 9).Value = 7296.5884
$ws.Cells.Item(134, 10).Value = 1000
$ws.Cells.Item(134, 11).Value = 21889.7652
$ws.Cells.Item(134, 12).Value = 3000
$ws.Cells.Item(134, 13).Value = -19354.7652
$ws.Cells.Item(134, 14).Value = -8070

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 77231.16
$ws.Cells.Item(38, 9).Value = 76.28570999999999
$ws.Cells.Item(38, 10).Value = 167245.17
$ws.Cells.Item(38, 11).Value = 228.85713
$ws.Cells.Item(38, 12).Value = 501735.51
$ws.Cells.Item(38, 13).Value = 118.14287
$ws.Cells.Item(38, 14).Value = -502429.51
$ws.Cells.Item(68, 8).Value = 806.38
$ws.Cells.Item(68, 9).Value = 644.2549
$ws.Cells.Item(68, 10).Value = 975.12244
$ws.Cells.Item(68, 11).Value = 1932.7647
$ws.Cells.Item(68, 12).Value = 2925.36732
$ws.Cells.Item(68, 13).Value = -1121.7647
$ws.Cells.Item(68, 14).Value = -4547.367319999999
$ws.Cells.Item(71, 8).Value = 806.38
$ws.Cells.Item(71, 9).Value = 644.2549
$ws.Cells.Item(71, 10).Value = 975.12244
$ws.Cells.Item(71, 11).Value = 5798.2941
$ws.Cells.Item(71, 12).Value = 8776.10196
$ws.Cells.Item(71, 13).Value = -1742.2941
$ws.Cells.Item(71, 14).Value = -16888.10196
$ws.Cells.Item(107, 8).Value = 936.5714
$ws.Cells.Item(107, 9).Value = 988.35297
$ws.Cells.Item(107, 10).Value = 887.6667
$ws.Cells.Item(107, 11).Value = 2965.05891
$ws.Cells.Item(107, 12).Value = 2663.0001
$ws.Cells.Item(107, 13).Value = -1045.05891
$ws.Cells.Item(107, 14).Value = -6503.0001
$ws.Cells.Item(131, 8).Value = 3337.3845
$ws.Cells.Item(131, 10).Value = 4082.805
$ws.Cells.Item(131, 12).Value = 12248.415
$ws.Cells.Item(131, 14).Value = -22328.415

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6237.7905
$ws.Cells.Item(70, 9).Value = 7110.871
$ws.Cells.Item(70, 10).Value = 3982.3333
$ws.Cells.Item(70, 11).Value = 7110.871
$ws.Cells.Item(70, 12).Value = 3982.3333
$ws.Cells.Item(70, 13).Value = -6840.871
$ws.Cells.Item(70, 14).Value = -4522.3333
$ws.Cells.Item(73, 8).Value = 6237.7905
$ws.Cells.Item(73, 9).Value = 7110.871
$ws.Cells.Item(73, 10).Value = 3982.3333
$ws.Cells.Item(73, 11).Value = 7110.871
$ws.Cells.Item(73, 12).Value = 3982.3333
$ws.Cells.Item(73, 13).Value = -6174.871
$ws.Cells.Item(73, 14).Value = -5854.3333
$ws.Cells.Item(132, 8).Value = 4172.364
$ws.Cells.Item(132, 9).Value = 3229.5557
$ws.Cells.Item(132, 10).Value = 4825.077
$ws.Cells.Item(132, 11).Value = 9688.667099999999
$ws.Cells.Item(132, 12).Value = 14475.231
$ws.Cells.Item(132, 13).Value = -7158.667099999999
$ws.Cells.Item(132, 14).Value = -19535.231
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(6, 8).Value = 28666.666
$ws.Cells.Item(6, 10).Value = 28666.666
$ws.Cells.Item(6, 12).Value = 28666.666
$ws.Cells.Item(6, 14).Value = -28890.666
$ws.Cells.Item(7, 8).Value = 40470
$ws.Cells.Item(7, 9).Value = 51075.285
$ws.Cells.Item(7, 10).Value = 3351.5
$ws.Cells.Item(7, 11).Value = 51075.285
$ws.Cells.Item(7, 12).Value = 3351.5
$ws.Cells.Item(7, 13).Value = -50963.285
$ws.Cells.Item(7, 14).Value = -3575.5
$ws.Cells.Item(21, 8).Value = 13105.286
$ws.Cells.Item(21, 9).Value = 346
$ws.Cells.Item(21, 10).Value = 45003.5
$ws.Cells.Item(21, 11).Value = 346
$ws.Cells.Item(21, 12).Value = 45003.5
$ws.Cells.Item(21, 13).Value = -172
$ws.Cells.Item(21, 14).Value = -45351.5
$ws.Cells.Item(61, 8).Value = 2347.25
$ws.Cells.Item(61, 9).Value = 1991.7
$ws.Cells.Item(61, 10).Value = 4125
$ws.Cells.Item(61, 11).Value = 1991.7
$ws.Cells.Item(61, 12).Value = 4125
$ws.Cells.Item(61, 13).Value = -1789.7
$ws.Cells.Item(61, 14).Value = -4529
$ws.Cells.Item(68, 8).Value = 1864
$ws.Cells.Item(68, 9).Value = 1872.4
$ws.Cells.Item(68, 10).Value = 1850
$ws.Cells.Item(68, 11).Value = 1872.4
$ws.Cells.Item(68, 12).Value = 1850
$ws.Cells.Item(68, 13).Value = -1123.4
$ws.Cells.Item(68, 14).Value = -3348
$ws.Cells.Item(71, 8).Value = 1864
$ws.Cells.Item(71, 9).Value = 1872.4
$ws.Cells.Item(71, 10).Value = 1850
$ws.Cells.Item(71, 11).Value = 9362
$ws.Cells.Item(71, 12).Value = 9250
$ws.Cells.Item(71, 13).Value = -5618
$ws.Cells.Item(71, 14).Value = -16738
$ws.Cells.Item(113, 8).Value = 2347.25
$ws.Cells.Item(113, 9).Value = 1991.7
$ws.Cells.Item(113, 10).Value = 4125
$ws.Cells.Item(113, 11).Value = 1991.7
$ws.Cells.Item(113, 12).Value = 4125
$ws.Cells.Item(113, 13).Value = 178.3
$ws.Cells.Item(113, 14).Value = -8465
$ws.Cells.Item(126, 8).Value = 40470
$ws.Cells.Item(126, 9).Value = 51075.285
$ws.Cells.Item(126, 10).Value = 3351.5
$ws.Cells.Item(126, 11).Value = 153225.855
$ws.Cells.Item(126, 12).Value = 10054.5
$ws.Cells.Item(126, 13).Value = -150755.855
$ws.Cells.Item(126, 14).Value = -14994.5
$ws.Cells.Item(132, 8).Value = 4030.2554
$ws.Cells.Item(132, 9).Value = 3495.7354
$ws.Cells.Item(132, 10).Value = 5428.231
$ws.Cells.Item(132, 11).Value = 10487.2062
$ws.Cells.Item(132, 12).Value = 16284.693
$ws.Cells.Item(132, 13).Value = -7957.206200000001
$ws.Cells.Item(132, 14).Value = -21344.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1719.4445
$ws.Cells.Item(113, 9).Value = 1834.375
$ws.Cells.Item(113, 11).Value = 5503.125
$ws.Cells.Item(113, 13).Value = -3333.125
$ws.Cells.Item(126, 8).Value = 36142.484
$ws.Cells.Item(126, 9).Value = 48820.668
$ws.Cells.Item(126, 10).Value = 2862.25
$ws.Cells.Item(126, 11).Value = 146462.004
$ws.Cells.Item(126, 12).Value = 8586.75
$ws.Cells.Item(126, 13).Value = -143992.004
$ws.Cells.Item(126, 14).Value = -13526.75
$ws.Cells.Item(132, 8).Value = 2092.963
$ws.Cells.Item(132, 9).Value = 1351.3871
$ws.Cells.Item(132, 10).Value = 3092.4783
$ws.Cells.Item(132, 11).Value = 4054.1613
$ws.Cells.Item(132, 12).Value = 9277.4349
$ws.Cells.Item(132, 13).Value = -1524.1613
$ws.Cells.Item(132, 14).Value = -14337.4349
$ws.Cells.Item(136, 8).Value = 5955670.5
$ws.Cells.Item(136, 9).Value = 13889380
$ws.Cells.Item(136, 10).Value = 5389.0625
$ws.Cells.Item(136, 11).Value = 41668140
$ws.Cells.Item(136, 12).Value = 16167.1875
$ws.Cells.Item(136, 13).Value = -41665590
$ws.Cells.Item(136, 14).Value = -21267.1875
